$wb = $excel.ActiveWorkbook

# --- EMP sheet: replace numeric DATA_LENGTH values with Oracle data-type/length strings ---
$wsEmp = $wb.Worksheets.Item("EMP")

$wsEmp.Range("C2").Value = "NUMBER(4)"
$wsEmp.Range("C3").Value = "VARCHAR2(10) "
$wsEmp.Range("C4").Value = "VARCHAR2(9) "
$wsEmp.Range("C5").Value = "NUMBER(4)"
$wsEmp.Range("C6").Value = "DATE"
$wsEmp.Range("C7").Value = "NUMBER(7,2)"
$wsEmp.Range("C8").Value = "NUMBER(7,2)"
$wsEmp.Range("C9").Value = "NUMBER(2)"

# --- DEPT sheet: same fix for Mysql/Oracle dept metadata ---
$wsDept = $wb.Worksheets.Item("DEPT")

$wsDept.Range("C2").Value = "NUMBER(2) "
$wsDept.Range("C3").Value = "VARCHAR2(14)"
$wsDept.Range("C4").Value = "VARCHAR2(13)"

# --- Selections / active tab: EMP was selected before, DEPT is selected now ---
$wsEmp.Activate()
$wsEmp.Range("C2:C9").Select()

$wsDept.Activate()
$wsDept.Range("C2:C4").Select()
